$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 246 (pushes existing rows 246..278 down to 247..279)
$ws.Rows.Item(246).Insert()

# Populate the newly inserted row 246 with the new record
$ws.Range("A246").Value = 4
$ws.Range("B246").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C246").Value = "Los Lagos"
$ws.Range("D246").Value = 44776
$ws.Range("E246").Value = 10
$ws.Range("F246").Value = "Fruta"
$ws.Range("G246").Value = 100108
$ws.Range("H246").Value = "Tropicales y subtropicales"
$ws.Range("I246").Value = 100108005
$ws.Range("J246").Value = "Piña"
$ws.Range("K246").Value = "Caramelo"
$ws.Range("L246").Value = "Primera"
$ws.Range("M246").Value = 30
$ws.Range("N246").Value = 23000
$ws.Range("O246").Value = 23000
$ws.Range("P246").Value = 23000
$ws.Range("Q246").Value = "$/caja 12 unidades"
$ws.Range("R246").Value = "Ecuador"
$ws.Range("S246").Value = 1917
$ws.Range("T246").Value = 12
